$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 57.433
$ws.Range("D2").Value = 57.433
$ws.Range("E2").Value = 2.99337846
$ws.Range("F2").Value = 0.01994631
$ws.Range("G2").Value = 1.13410404
$ws.Range("H2").Value = 65.52918533
$ws.Range("I2").Value = 7.424327015402246
$ws.Range("J2").Value = 7.424327015402246
$ws.Range("K2").Value = 0.4022620362848879
$ws.Range("L2").Value = 0.002377418011267749
$ws.Range("M2").Value = 0.1201558446513622
$ws.Range("N2").Value = 12.70321905575497
$ws.Range("C3").Value = 90.00700000000001
$ws.Range("D3").Value = 90.00700000000001
$ws.Range("E3").Value = 1.9146389
$ws.Range("F3").Value = 0.00759413
$ws.Range("G3").Value = 0.6693451199999999
$ws.Range("H3").Value = 60.16411374
$ws.Range("I3").Value = 12.46171116842487
$ws.Range("J3").Value = 12.46171116842487
$ws.Range("K3").Value = 0.2734552175076793
$ws.Range("L3").Value = 0.001578775086538255
$ws.Range("M3").Value = 0.09240722003230969
$ws.Range("N3").Value = 11.35750868439249
$ws.Range("C4").Value = 27.766
$ws.Range("D4").Value = 55.508
$ws.Range("E4").Value = 3.176143
$ws.Range("F4").Value = 0.03767816
$ws.Range("G4").Value = 0.51135439
$ws.Range("H4").Value = 14.47240166
$ws.Range("I4").Value = 5.787463929800231
$ws.Range("J4").Value = 11.5715096610527
$ws.Range("K4").Value = 0.6555328995299673
$ws.Range("L4").Value = 0.005634863361209347
$ws.Range("M4").Value = 0.07158144819821721
$ws.Range("N4").Value = 4.634190969150257
$ws.Range("C5").Value = 48.055
$ws.Range("D5").Value = 94.18899999999999
$ws.Range("E5").Value = 1.8370735
$ws.Range("F5").Value = 0.01204469
$ws.Range("G5").Value = 0.28201036
$ws.Range("H5").Value = 13.55966388
$ws.Range("I5").Value = 7.903191471515506
$ws.Range("J5").Value = 14.53301370973989
$ws.Range("K5").Value = 0.2839911382715432
$ws.Range("L5").Value = 0.002551167803391432
$ws.Range("M5").Value = 0.0384626357549716
$ws.Range("N5").Value = 2.963306803831891
$ws.Range("C6").Value = 14.489
$ws.Range("D6").Value = 57.899
$ws.Range("E6").Value = 3.12849664
$ws.Range("F6").Value = 0.06343076
$ws.Range("G6").Value = 0.22269545
$ws.Range("H6").Value = 3.35020085
$ws.Range("I6").Value = 3.851585869964957
$ws.Range("J6").Value = 15.38946505355643
$ws.Range("K6").Value = 0.827824253644596
$ws.Range("L6").Value = 0.009974247261364179
$ws.Range("M6").Value = 0.03998941746145356
$ws.Range("N6").Value = 1.446412518733341
$ws.Range("C7").Value = 25.459
$ws.Range("D7").Value = 93.842
$ws.Range("E7").Value = 1.86178314
$ws.Range("F7").Value = 0.01908867
$ws.Range("G7").Value = 0.11673369
$ws.Range("H7").Value = 2.99794069
$ws.Range("I7").Value = 5.876976931880697
$ws.Range("J7").Value = 16.99835227249924
$ws.Range("K7").Value = 0.3478086949497137
$ws.Range("L7").Value = 0.005231700509520908
$ws.Range("M7").Value = 0.0240389042554294
$ws.Range("N7").Value = 1.057252794355634
$ws.Range("C8").Value = 9.351000000000001
$ws.Range("D8").Value = 56.026
$ws.Range("E8").Value = 3.20859392
$ws.Range("F8").Value = 0.09406307999999999
$ws.Range("G8").Value = 0.14364999
$ws.Range("H8").Value = 1.40160228
$ws.Range("I8").Value = 2.436336334485601
$ws.Range("J8").Value = 14.60652849498041
$ws.Range("K8").Value = 0.7848792580316554
$ws.Range("L8").Value = 0.01451697287367497
$ws.Range("M8").Value = 0.03151001155391771
$ws.Range("N8").Value = 0.6633873860143291
$ws.Range("C9").Value = 17.314
$ws.Range("D9").Value = 88.161
$ws.Range("E9").Value = 1.9772961
$ws.Range("F9").Value = 0.021473
$ws.Range("G9").Value = 0.05972508
$ws.Range("H9").Value = 1.05217963
$ws.Range("I9").Value = 4.177182405982872
$ws.Range("J9").Value = 15.50411095343607
$ws.Range("K9").Value = 0.3532137840572449
$ws.Range("L9").Value = 0.005968395105900759
$ws.Range("M9").Value = 0.01437331033699765
$ws.Range("N9").Value = 0.4332810593949307
$ws.Range("C10").Value = 6.812
$ws.Range("D10").Value = 54.362
$ws.Range("E10").Value = 3.30291711
$ws.Range("F10").Value = 0.11096264
$ws.Range("G10").Value = 0.0934958
$ws.Range("H10").Value = 0.66820472
$ws.Range("I10").Value = 1.720353340339018
$ws.Range("J10").Value = 13.74954566543949
$ws.Range("K10").Value = 0.8189698617737049
$ws.Range("L10").Value = 0.0131133136096192
$ws.Range("M10").Value = 0.02169153087432216
$ws.Range("N10").Value = 0.3163854945162023
$ws.Range("C11").Value = 13.272
$ws.Range("D11").Value = 81.985
$ws.Range("E11").Value = 2.13745123
$ws.Range("F11").Value = 0.02152434
$ws.Range("G11").Value = 0.03429336000000001
$ws.Range("H11").Value = 0.46782531
$ws.Range("I11").Value = 3.575024027443284
$ws.Range("J11").Value = 15.65118956292891
$ws.Range("K11").Value = 0.413616642121212
$ws.Range("L11").Value = 0.007111314545072448
$ws.Range("M11").Value = 0.01086492314124979
$ws.Range("N11").Value = 0.2375206677700318
$ws.Range("C12").Value = 5.677
$ws.Range("D12").Value = 56.619
$ws.Range("E12").Value = 3.18791617
$ws.Range("F12").Value = 0.1287425
$ws.Range("G12").Value = 0.07289238999999999
$ws.Range("H12").Value = 0.43940232
$ws.Range("I12").Value = 1.476770751285335
$ws.Range("J12").Value = 14.75573009276635
$ws.Range("K12").Value = 0.8301048504131936
$ws.Range("L12").Value = 0.01459983195256175
$ws.Range("M12").Value = 0.01968038161745113
$ws.Range("N12").Value = 0.2288489329110191
$ws.Range("C13").Value = 10.875
$ws.Range("D13").Value = 75.901
$ws.Range("E13").Value = 2.31650139
$ws.Range("F13").Value = 0.02081715
$ws.Range("G13").Value = 0.02178227
$ws.Range("H13").Value = 0.2474617
$ws.Range("I13").Value = 3.285143401309472
$ws.Range("J13").Value = 15.46616957754574
$ws.Range("K13").Value = 0.4579194641580796
$ws.Range("L13").Value = 0.007711707476928162
$ws.Range("M13").Value = 0.008324527680788155
$ws.Range("N13").Value = 0.1503671309960857
